$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 632, shifting existing rows 632:693 down to 633:694
$ws.Rows.Item(632).Insert()

# Populate the newly inserted row 632 with the new record's data
$ws.Range("A632").Value = 6
$ws.Range("B632").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C632").Value = 'Metropolitana'
$ws.Range("D632").Value = 45223
$ws.Range("E632").Value = 13
$ws.Range("F632").Value = 100112043
$ws.Range("G632").Value = 'Pepino ensalada'
$ws.Range("H632").Value = 'Sin especificar'
$ws.Range("I632").Value = 'Primera'
$ws.Range("J632").Value = 360
$ws.Range("K632").Value = 11000
$ws.Range("L632").Value = 12000
$ws.Range("M632").Value = 11333
$ws.Range("N632").Value = '$/caja 60 unidades'
$ws.Range("O632").Value = 'Región de Arica y Parinacota'
$ws.Range("P632").Value = 189
$ws.Range("Q632").Value = 60
$ws.Range("R632").Value = 'Hortaliza'
